$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A87").Value = "LIVE, SEA, BEACH"
$ws.Range("B87").Value = "-33.01330799002186, -71.55553066972728"
$ws.Range("C87").Value = "Playa Acapulco"
$ws.Range("D87").Value = "Viña del Mar"
$ws.Range("E87").Value = "Chile"
$ws.Range("F87").Value = "bmSBriUdJzc"

$ws.Range("A87").Style = $ws.Range("A86").Style
$ws.Range("C87:E87").Style = $ws.Range("C86:E86").Style
